$wb = $excel.ActiveWorkbook

# Add the new, empty "2012" multiform sheet at the end of the workbook
# (handles the previously-unhandled "empty multiform" case) and populate
# its single informational cell.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "2012"
$newSheet.Range("A1").Value = "aucune instruction"

# Make the newly added sheet the active one, matching where the cursor
# was left on that sheet.
$newSheet.Activate()
$newSheet.Range("D11").Select() | Out-Null
